$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Suite")

# Update C6 and C7 from "N" to "Y"
$ws.Range("C6").Value = "Y"
$ws.Range("C7").Value = "Y"

# Update the selection to C2:C7 with active cell C2
$ws.Range("C2:C7").Select()
